{"js": "// Update the date and the multiplication problems to the new values.\nconst replacements = [\n  [\"2024-12-13 Friday\", \"2024-12-14 Saturday\"],\n  [\"411\u00d72=\", \"266\u00d77=\"],\n  [\"409\u00d72=\", \"423\u00d75=\"],\n  [\"625\u00d72=\", \"927\u00d74=\"],\n  [\"288\u00d72=\", \"729\u00d78=\"],\n  [\"862\u00d76=\", \"788\u00d79=\"],\n  [\"231\u00d79=\", \"561\u00d76=\"],\n  [\"266\u00d76=\", \"820\u00d75=\"],\n  [\"149\u00d74=\", \"105\u00d78=\"],\n  [\"292\u00d76=\", \"105\u00d72=\"],\n  [\"644\u00d74=\", \"386\u00d75=\"],\n  [\"711\u00d78=\", \"933\u00d74=\"],\n  [\"462\u00d77=\", \"297\u00d73=\"],\n  [\"225\u00d78=\", \"870\u00d72=\"],\n  [\"995\u00d72=\", \"947\u00d72=\"],\n  [\"826\u00d78=\", \"341\u00d77=\"],\n  [\"237\u00d76=\", \"444\u00d79=\"],\n  [\"332\u00d77=\", \"629\u00d73=\"],\n  [\"955\u00d74=\", \"755\u00d78=\"],\n  [\"338\u00d72=\", \"688\u00d73=\"],\n  [\"123\u00d75=\", \"583\u00d72=\"],\n  [\"276\u00d74=\", \"863\u00d72=\"],\n  [\"345\u00d73=\", \"713\u00d75=\"],\n  [\"576\u00d73=\", \"192\u00d77=\"],\n  [\"736\u00d75=\", \"275\u00d73=\"],\n  [\"160\u00d73=\", \"531\u00d79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date and the multiplication problems to the new values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-12-13 Friday\", \"2024-12-14 Saturday\"),\n    @(\"411\u00d72=\", \"266\u00d77=\"),\n    @(\"409\u00d72=\", \"423\u00d75=\"),\n    @(\"625\u00d72=\", \"927\u00d74=\"),\n    @(\"288\u00d72=\", \"729\u00d78=\"),\n    @(\"862\u00d76=\", \"788\u00d79=\"),\n    @(\"231\u00d79=\", \"561\u00d76=\"),\n    @(\"266\u00d76=\", \"820\u00d75=\"),\n    @(\"149\u00d74=\", \"105\u00d78=\"),\n    @(\"292\u00d76=\", \"105\u00d72=\"),\n    @(\"644\u00d74=\", \"386\u00d75=\"),\n    @(\"711\u00d78=\", \"933\u00d74=\"),\n    @(\"462\u00d77=\", \"297\u00d73=\"),\n    @(\"225\u00d78=\", \"870\u00d72=\"),\n    @(\"995\u00d72=\", \"947\u00d72=\"),\n    @(\"826\u00d78=\", \"341\u00d77=\"),\n    @(\"237\u00d76=\", \"444\u00d79=\"),\n    @(\"332\u00d77=\", \"629\u00d73=\"),\n    @(\"955\u00d74=\", \"755\u00d78=\"),\n    @(\"338\u00d72=\", \"688\u00d73=\"),\n    @(\"123\u00d75=\", \"583\u00d72=\"),\n    @(\"276\u00d74=\", \"863\u00d72=\"),\n    @(\"345\u00d73=\", \"713\u00d75=\"),\n    @(\"576\u00d73=\", \"192\u00d77=\"),\n    @(\"736\u00d75=\", \"275\u00d73=\"),\n    @(\"160\u00d73=\", \"531\u00d79=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
